# Updates cryptos list data (Price / Volume(1h) columns, and a Stellar/OKB row swap)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "89.728.41"
$ws.Cells.Item(2, 5).Value = "  +3.06%  "
$ws.Cells.Item(3, 4).Value = "3.299.86"
$ws.Cells.Item(3, 5).Value = "  -1.77%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "213.15"
$ws.Cells.Item(5, 5).Value = "  -3.54%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "629.83"
$ws.Cells.Item(6, 5).Value = "  -1.52%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.384"
$ws.Cells.Item(7, 5).Value = "  +17.78%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.717"
$ws.Cells.Item(8, 5).Value = "  +14.55%  "
$ws.Cells.Item(9, 5).Value = "  +0.06%  "
$ws.Cells.Item(10, 4).Value = "3.295.91"
$ws.Cells.Item(10, 5).Value = "  -2.45%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.581"
$ws.Cells.Item(11, 5).Value = "  -5.73%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.187"
$ws.Cells.Item(12, 5).Value = "  +11.94%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000265"
$ws.Cells.Item(13, 5).Value = "  -3.97%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "34.38"
$ws.Cells.Item(14, 5).Value = "  -0.59%  "
$ws.Cells.Item(15, 4).Value = "3.896.42"
$ws.Cells.Item(15, 5).Value = "  -1.84%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.46"
$ws.Cells.Item(16, 5).Value = "  +0.59%  "
$ws.Cells.Item(17, 4).Value = "89.337.84"
$ws.Cells.Item(17, 5).Value = "  +2.83%  "
$ws.Cells.Item(18, 4).Value = "3.303.49"
$ws.Cells.Item(18, 5).Value = "  -1.57%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.17"
$ws.Cells.Item(19, 5).Value = "  -4.26%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "3.09"
$ws.Cells.Item(20, 5).Value = "  -4.18%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "439.04"
$ws.Cells.Item(21, 5).Value = "  -2.30%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "8.96"
$ws.Cells.Item(22, 5).Value = "  -2.97%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.37"
$ws.Cells.Item(23, 5).Value = "  +0.43%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "7.45"
$ws.Cells.Item(24, 5).Value = "  +0.14%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "5.28"
$ws.Cells.Item(25, 5).Value = "  -3.03%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "12.25"
$ws.Cells.Item(26, 5).Value = "  -1.49%  "
$ws.Cells.Item(27, 5).Value = "  +0.36%  "
$ws.Cells.Item(28, 5).Value = "  +4.00%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "77.27"
$ws.Cells.Item(29, 5).Value = "  -1.95%  "
$ws.Cells.Item(30, 5).Value = "  -0.03%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.184"
$ws.Cells.Item(31, 5).Value = "  -0.63%  "
$ws.Cells.Item(32, 5).Value = "  +0.02%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "8.91"
$ws.Cells.Item(33, 5).Value = "  -4.86%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "564.24"
$ws.Cells.Item(34, 5).Value = "  -7.22%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.39"
$ws.Cells.Item(35, 5).Value = "  -11.85%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "7.18"
$ws.Cells.Item(36, 5).Value = "  +8.96%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.97"
$ws.Cells.Item(37, 5).Value = "  -4.63%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.141"
$ws.Cells.Item(38, 5).Value = "  -7.56%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "22.81"
$ws.Cells.Item(39, 5).Value = "  -3.40%  "
$ws.Cells.Item(40, 5).Value = "  +2.40%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.999"
$ws.Cells.Item(41, 5).Value = "  +0.15%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.13"
$ws.Cells.Item(42, 5).Value = "  -0.63%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.404"
$ws.Cells.Item(43, 5).Value = "  -4.75%  "
$ws.Cells.Item(44, 5).Value = "  -2.72%  "
$ws.Cells.Item(45, 5).Value = "  -0.01%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "153.66"
$ws.Cells.Item(46, 5).Value = "  -2.12%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "182.19"
$ws.Cells.Item(47, 5).Value = "  -4.79%  "
$ws.Cells.Item(48, 2).Value = "OKB"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "45.03"
$ws.Cells.Item(48, 5).Value = "  -1.55%  "
$ws.Cells.Item(49, 2).Value = "Stellar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.132"
$ws.Cells.Item(49, 5).Value = "  +16.46%  "
$ws.Cells.Item(50, 5).Value = "  -4.57%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.25"
$ws.Cells.Item(51, 5).Value = "  -1.70%  "
